$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "чулок"
$ws.Range("A25").Value = "выбойка"
$ws.Range("A26").Value = "сахар"
$ws.Range("A27").Value = "шелк"
$ws.Range("A28").Value = "лыко"
$ws.Range("A29").Value = "лес"
$ws.Range("A30").Value = "ладан"
$ws.Range("A31").Value = "сапог"
$ws.Range("A32").Value = "китайка"
$ws.Range("A33").Value = "коса"
$ws.Range("A34").Value = "сани"
$ws.Range("A35").Value = "платок"
$ws.Range("A36").Value = "ром"
$ws.Range("A37").Value = "обод"
$ws.Range("A38").Value = "конь"
$ws.Range("A39").Value = "гвоздь"
$ws.Range("A40").Value = "веревка"
$ws.Range("A41").Value = "овца"
$ws.Range("A42").Value = "горшок"
$ws.Range("A43").Value = "рогожа"
$ws.Range("A44").Value = "замок"
$ws.Range("A45").Value = "сосуд"
$ws.Range("A46").Value = "покроми"
$ws.Range("A47").Value = "скотский кожа"
$ws.Range("A48").Value = "гумми"
$ws.Range("A49").Value = "нитка"
$ws.Range("A50").Value = "сковорода"
$ws.Range("A51").Value = "дуга"
$ws.Range("A52").Value = "котел"
$ws.Range("A53").Value = "роза"
$ws.Range("A54").Value = "хомут"
$ws.Range("A55").Value = "бечева"
$ws.Range("A56").Value = "брусья"
